$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the theta_threshold_range row (row 5) entirely - this shifts
# pie_threshold_range (old row 6) up to row 5.
$ws.Rows.Item(5).Delete()

# Update the remaining Min/Max values that changed.
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 4.0999999999999996
$ws.Range("B4").Value = 0.7
$ws.Range("C4").Value = 1.4

# Match the saved selection state (B4:C4 active cell B4).
$ws.Range("B4:C4").Select()

# Page setup now present in the sheet (paper size / orientation / DPI).
$ps = $ws.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1
